$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.969.36"
$ws.Range("E2").Value = "'  +1.59%  "
$ws.Range("D3").Value = "'3.355.69"
$ws.Range("E3").Value = "'  +0.89%  "
$ws.Range("E4").Value = "'  +0.06%  "
$ws.Range("D5").Value = "'586.69"
$ws.Range("E5").Value = "'  +1.77%  "
$ws.Range("D6").Value = "'178.03"
$ws.Range("E6").Value = "'  +1.53%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "'  +0.00%  "
$ws.Range("D8").Value = "'0.594"
$ws.Range("E8").Value = "'  +1.06%  "
$ws.Range("D9").Value = "'0.192"
$ws.Range("E9").Value = "'  +7.47%  "
$ws.Range("D10").Value = "'0.585"
$ws.Range("E10").Value = "'  +1.43%  "
$ws.Range("D11").Value = "'48.29"
$ws.Range("E11").Value = "'  +5.15%  "
$ws.Range("D12").Value = "'0.0000279"
$ws.Range("E12").Value = "'  +3.41%  "
$ws.Range("D13").Value = "'701.39"
$ws.Range("E13").Value = "'  +3.00%  "
$ws.Range("B14").Value = "'Polkadot"
$ws.Range("C14").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'8.52"
$ws.Range("E14").Value = "'  +1.53%  "
$ws.Range("B15").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "'3.898.79"
$ws.Range("E15").Value = "'  +0.74%  "
$ws.Range("D16").Value = "'68.972.92"
$ws.Range("E17").Value = "'  +1.59%  "
$ws.Range("D18").Value = "'3.361.16"
$ws.Range("E18").Value = "'  +1.10%  "
$ws.Range("D19").Value = "'17.59"
$ws.Range("E19").Value = "'  +1.41%  "
$ws.Range("D20").Value = "'11.36"
$ws.Range("E20").Value = "'  +3.88%  "
$ws.Range("D21").Value = "'0.898"
$ws.Range("E21").Value = "'  +1.17%  "
$ws.Range("D22").Value = "'5.45"
$ws.Range("E22").Value = "'  +2.10%  "
$ws.Range("D23").Value = "'17.03"
$ws.Range("E23").Value = "'  -0.03%  "
$ws.Range("D24").Value = "'101.20"
$ws.Range("E24").Value = "'  +2.81%  "
$ws.Range("D25").Value = "'3.94"
$ws.Range("E25").Value = "'  +2.38%  "
$ws.Range("E26").Value = "'  +1.02%  "
$ws.Range("D27").Value = "'9.58"
$ws.Range("E27").Value = "'  +2.92%  "
$ws.Range("D28").Value = "'33.63"
$ws.Range("E28").Value = "'  +1.83%  "
$ws.Range("D29").Value = "'8.63"
$ws.Range("E29").Value = "'  +1.89%  "
$ws.Range("D30").Value = "'7.06"
$ws.Range("E30").Value = "'  -1.74%  "
$ws.Range("D31").Value = "'11.12"
$ws.Range("E31").Value = "'  +1.72%  "
$ws.Range("D32").Value = "'552.76"
$ws.Range("E32").Value = "'  -3.19%  "
$ws.Range("E33").Value = "'  +0.83%  "
$ws.Range("D34").Value = "'3.51"
$ws.Range("E34").Value = "'  +7.29%  "
$ws.Range("D35").Value = "'57.49"
$ws.Range("E35").Value = "'  +0.03%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "'  -0.04%  "
$ws.Range("D37").Value = "'3.697.87"
$ws.Range("E37").Value = "'  -0.71%  "
$ws.Range("D38").Value = "'0.141"
$ws.Range("E38").Value = "'  +7.70%  "
$ws.Range("D39").Value = "'35.06"
$ws.Range("E39").Value = "'  +2.15%  "
$ws.Range("E40").Value = "'  +1.05%  "
$ws.Range("D41").Value = "'2.65"
$ws.Range("E41").Value = "'  -0.03%  "
$ws.Range("D42").Value = "'0.0₃0689"
$ws.Range("E42").Value = "'  +3.05%  "
$ws.Range("E43").Value = "'  +1.08%  "
$ws.Range("D44").Value = "'0.0418"
$ws.Range("E44").Value = "'  +2.99%  "
$ws.Range("D45").Value = "'3.27"
$ws.Range("E45").Value = "'  -2.08%  "
$ws.Range("E46").Value = "'  +0.74%  "
$ws.Range("E47").Value = "'  +0.99%  "
$ws.Range("E48").Value = "'  -0.05%  "
$ws.Range("D49").Value = "'1.36"
$ws.Range("E49").Value = "'  +1.22%  "
$ws.Range("D50").Value = "'132.83"
$ws.Range("E50").Value = "'  +4.10%  "
$ws.Range("D51").Value = "'2.62"
$ws.Range("E51").Value = "'  -3.74%  "
